{"js": "// Replace each table-cell expression with its updated result, in document order.\n// Table is 20 rows x 5 columns = 100 cells; `replacements[i]` is [oldText, newText]\n// for the cell at row = floor(i/5), col = i%5 (0-based).\nconst replacements = [\n  [\"86-71=15\", \"36+59=95\"],\n  [\"38-31=7\", \"90-54=36\"],\n  [\"60-17=43\", \"99-2=97\"],\n  [\"22+28=50\", \"66-7=59\"],\n  [\"74-40=34\", \"69-14=55\"],\n  [\"61-21=40\", \"59+35=94\"],\n  [\"98-84=14\", \"35+56=91\"],\n  [\"18+62=80\", \"60+11=71\"],\n  [\"84+12=96\", \"23+11=34\"],\n  [\"59+17=76\", \"44-22=22\"],\n  [\"28+11=39\", \"55+39=94\"],\n  [\"78-62=16\", \"10+78=88\"],\n  [\"73-67=6\", \"76+6=82\"],\n  [\"87+2=89\", \"59-11=48\"],\n  [\"91-67=24\", \"59-5=54\"],\n  [\"97-89=8\", \"28+14=42\"],\n  [\"46+30=76\", \"3+58=61\"],\n  [\"58-51=7\", \"42+2=44\"],\n  [\"85+8=93\", \"9+42=51\"],\n  [\"3+75=78\", \"80-58=22\"],\n  [\"66-29=37\", \"55+6=61\"],\n  [\"53-39=14\", \"25+18=43\"],\n  [\"13+85=98\", \"7+85=92\"],\n  [\"34+38=72\", \"15+55=70\"],\n  [\"5+75=80\", \"7+33=40\"],\n  [\"17+25=42\", \"53-3=50\"],\n  [\"39+3=42\", \"68-65=3\"],\n  [\"77-23=54\", \"90-53=37\"],\n  [\"33+48=81\", \"76-49=27\"],\n  [\"98-46=52\", \"45-43=2\"],\n  [\"96-35=61\", \"31-0=31\"],\n  [\"56+29=85\", \"70-10=60\"],\n  [\"71-64=7\", \"43+45=88\"],\n  [\"42-24=18\", \"48-41=7\"],\n  [\"14+61=75\", \"73-58=15\"],\n  [\"96-69=27\", \"82-69=13\"],\n  [\"49-0=49\", \"31-6=25\"],\n  [\"54+20=74\", \"8-5=3\"],\n  [\"18+28=46\", \"83-35=48\"],\n  [\"39+35=74\", \"47+17=64\"],\n  [\"73-70=3\", \"68-36=32\"],\n  [\"39-10=29\", \"42+55=97\"],\n  [\"72-7=65\", \"12+78=90\"],\n  [\"75+10=85\", \"17+13=30\"],\n  [\"37+5=42\", \"73-21=52\"],\n  [\"48-8=40\", \"16+79=95\"],\n  [\"93-18=75\", \"61+22=83\"],\n  [\"37+8=45\", \"47-40=7\"],\n  [\"91-14=77\", \"39+37=76\"],\n  [\"96-95=1\", \"93-21=72\"],\n  [\"92-89=3\", \"87-27=60\"],\n  [\"71+21=92\", \"84+5=89\"],\n  [\"34+55=89\", \"89-11=78\"],\n  [\"98-45=53\", \"67+28=95\"],\n  [\"93-52=41\", \"7+37=44\"],\n  [\"41+24=65\", \"69-50=19\"],\n  [\"89-67=22\", \"89-65=24\"],\n  [\"91-90=1\", \"16+73=89\"],\n  [\"11+59=70\", \"80-47=33\"],\n  [\"18+56=74\", \"6+91=97\"],\n  [\"31-11=20\", \"54+36=90\"],\n  [\"10+62=72\", \"43+21=64\"],\n  [\"6+75=81\", \"37-22=15\"],\n  [\"69-53=16\", \"65+27=92\"],\n  [\"42+38=80\", \"47+12=59\"],\n  [\"80-20=60\", \"79-3=76\"],\n  [\"1+51=52\", \"94-1=93\"],\n  [\"72+15=87\", \"18+22=40\"],\n  [\"95-9=86\", \"46-35=11\"],\n  [\"16+63=79\", \"22-20=2\"],\n  [\"3+73=76\", \"14+36=50\"],\n  [\"51+48=99\", \"35+27=62\"],\n  [\"21+3=24\", \"86+4=90\"],\n  [\"66-27=39\", \"28-6=22\"],\n  [\"69-28=41\", \"71-22=49\"],\n  [\"60-9=51\", \"39+15=54\"],\n  [\"60-9=51\", \"14+75=89\"],\n  [\"99-14=85\", \"38+27=65\"],\n  [\"34+8=42\", \"61+26=87\"],\n  [\"71-71=0\", \"40+58=98\"],\n  [\"38-7=31\", \"65+8=73\"],\n  [\"52-30=22\", \"46+15=61\"],\n  [\"58-14=44\", \"41+26=67\"],\n  [\"12+1=13\", \"56-47=9\"],\n  [\"49-21=28\", \"74+1=75\"],\n  [\"57+22=79\", \"99-53=46\"],\n  [\"1+2=3\", \"81-71=10\"],\n  [\"36+38=74\", \"95-8=87\"],\n  [\"7+42=49\", \"57-49=8\"],\n  [\"84-12=72\", \"33+54=87\"],\n  [\"30-23=7\", \"54+17=71\"],\n  [\"98-74=24\", \"27-2=25\"],\n  [\"52-35=17\", \"3+62=65\"],\n  [\"6+23=29\", \"40-10=30\"],\n  [\"84-13=71\", \"44-24=20\"],\n  [\"18+50=68\", \"64+6=70\"],\n  [\"42+34=76\", \"83+4=87\"],\n  [\"5+60=65\", \"0+92=92\"],\n  [\"50+2=52\", \"31+29=60\"],\n  [\"61+32=93\", \"52-23=29\"],\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\nconst table = tables.items[0];\n\nconst COLS = 5;\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const row = Math.floor(i / COLS);\n  const col = i % COLS;\n  const cell = table.getCell(row, col);\n  cell.body.paragraphs.load(\"items\");\n  await context.sync();\n  const para = cell.body.paragraphs.items[0];\n\n  // Search for the exact old expression text within the cell's paragraph and\n  // replace it in place (\"Replace\") so the run's existing formatting (font,\n  // size, paragraph alignment, etc.) is preserved - only the text changes.\n  const results = para.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Cell (${row},${col}) expected text \"${oldText}\" not found`);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace each table-cell expression with its updated result, in document order.\n# Table is 20 rows x 5 columns = 100 cells; each pair is (oldText, newText) for the\n# cell at the corresponding position (1-based row/col) when walking row-major.\n$replacements = @(\n    ,@('86-71=15', '36+59=95')\n    ,@('38-31=7', '90-54=36')\n    ,@('60-17=43', '99-2=97')\n    ,@('22+28=50', '66-7=59')\n    ,@('74-40=34', '69-14=55')\n    ,@('61-21=40', '59+35=94')\n    ,@('98-84=14', '35+56=91')\n    ,@('18+62=80', '60+11=71')\n    ,@('84+12=96', '23+11=34')\n    ,@('59+17=76', '44-22=22')\n    ,@('28+11=39', '55+39=94')\n    ,@('78-62=16', '10+78=88')\n    ,@('73-67=6', '76+6=82')\n    ,@('87+2=89', '59-11=48')\n    ,@('91-67=24', '59-5=54')\n    ,@('97-89=8', '28+14=42')\n    ,@('46+30=76', '3+58=61')\n    ,@('58-51=7', '42+2=44')\n    ,@('85+8=93', '9+42=51')\n    ,@('3+75=78', '80-58=22')\n    ,@('66-29=37', '55+6=61')\n    ,@('53-39=14', '25+18=43')\n    ,@('13+85=98', '7+85=92')\n    ,@('34+38=72', '15+55=70')\n    ,@('5+75=80', '7+33=40')\n    ,@('17+25=42', '53-3=50')\n    ,@('39+3=42', '68-65=3')\n    ,@('77-23=54', '90-53=37')\n    ,@('33+48=81', '76-49=27')\n    ,@('98-46=52', '45-43=2')\n    ,@('96-35=61', '31-0=31')\n    ,@('56+29=85', '70-10=60')\n    ,@('71-64=7', '43+45=88')\n    ,@('42-24=18', '48-41=7')\n    ,@('14+61=75', '73-58=15')\n    ,@('96-69=27', '82-69=13')\n    ,@('49-0=49', '31-6=25')\n    ,@('54+20=74', '8-5=3')\n    ,@('18+28=46', '83-35=48')\n    ,@('39+35=74', '47+17=64')\n    ,@('73-70=3', '68-36=32')\n    ,@('39-10=29', '42+55=97')\n    ,@('72-7=65', '12+78=90')\n    ,@('75+10=85', '17+13=30')\n    ,@('37+5=42', '73-21=52')\n    ,@('48-8=40', '16+79=95')\n    ,@('93-18=75', '61+22=83')\n    ,@('37+8=45', '47-40=7')\n    ,@('91-14=77', '39+37=76')\n    ,@('96-95=1', '93-21=72')\n    ,@('92-89=3', '87-27=60')\n    ,@('71+21=92', '84+5=89')\n    ,@('34+55=89', '89-11=78')\n    ,@('98-45=53', '67+28=95')\n    ,@('93-52=41', '7+37=44')\n    ,@('41+24=65', '69-50=19')\n    ,@('89-67=22', '89-65=24')\n    ,@('91-90=1', '16+73=89')\n    ,@('11+59=70', '80-47=33')\n    ,@('18+56=74', '6+91=97')\n    ,@('31-11=20', '54+36=90')\n    ,@('10+62=72', '43+21=64')\n    ,@('6+75=81', '37-22=15')\n    ,@('69-53=16', '65+27=92')\n    ,@('42+38=80', '47+12=59')\n    ,@('80-20=60', '79-3=76')\n    ,@('1+51=52', '94-1=93')\n    ,@('72+15=87', '18+22=40')\n    ,@('95-9=86', '46-35=11')\n    ,@('16+63=79', '22-20=2')\n    ,@('3+73=76', '14+36=50')\n    ,@('51+48=99', '35+27=62')\n    ,@('21+3=24', '86+4=90')\n    ,@('66-27=39', '28-6=22')\n    ,@('69-28=41', '71-22=49')\n    ,@('60-9=51', '39+15=54')\n    ,@('60-9=51', '14+75=89')\n    ,@('99-14=85', '38+27=65')\n    ,@('34+8=42', '61+26=87')\n    ,@('71-71=0', '40+58=98')\n    ,@('38-7=31', '65+8=73')\n    ,@('52-30=22', '46+15=61')\n    ,@('58-14=44', '41+26=67')\n    ,@('12+1=13', '56-47=9')\n    ,@('49-21=28', '74+1=75')\n    ,@('57+22=79', '99-53=46')\n    ,@('1+2=3', '81-71=10')\n    ,@('36+38=74', '95-8=87')\n    ,@('7+42=49', '57-49=8')\n    ,@('84-12=72', '33+54=87')\n    ,@('30-23=7', '54+17=71')\n    ,@('98-74=24', '27-2=25')\n    ,@('52-35=17', '3+62=65')\n    ,@('6+23=29', '40-10=30')\n    ,@('84-13=71', '44-24=20')\n    ,@('18+50=68', '64+6=70')\n    ,@('42+34=76', '83+4=87')\n    ,@('5+60=65', '0+92=92')\n    ,@('50+2=52', '31+29=60')\n    ,@('61+32=93', '52-23=29')\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$cols = $t.Columns.Count\n\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $row = [int][Math]::Floor($i / $cols) + 1\n    $col = ($i % $cols) + 1\n    $pair = $replacements[$i]\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $cell = $t.Cell($row, $col)\n    # Cell.Range.Text includes trailing cell-mark characters (CR + BEL); strip them\n    # before comparing so we can verify we're editing the expected cell.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $old) {\n        throw \"Cell ($row,$col) expected `\"$old`\" but found `\"$current`\"\"\n    }\n\n    # Assigning Range.Text replaces just the text and keeps the run's existing\n    # formatting (font, size) and the paragraph's properties (alignment) intact.\n    $cell.Range.Text = $new\n}"}
